$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: swap the "average_doctor" / "average_doctor_old" labels
$ws.Cells.Item(1, 68).Value = "average_doctor_old"
$ws.Cells.Item(1, 69).Value = "average_doctor"

$ws.Cells.Item(4, 5).Value = 0.5590000000000001
$ws.Cells.Item(4, 6).Value = 0.067
$ws.Cells.Item(4, 7).Value = 0.26
$ws.Cells.Item(4, 14).Value = 0.444
$ws.Cells.Item(4, 15).Value = 0.049
$ws.Cells.Item(4, 16).Value = 0.222
$ws.Cells.Item(4, 17).Value = 0.148
$ws.Cells.Item(4, 18).Value = 0.102
$ws.Cells.Item(4, 19).Value = 0.319
$ws.Cells.Item(4, 23).Value = 0.296
$ws.Cells.Item(4, 24).Value = 0.079
$ws.Cells.Item(4, 25).Value = 0.281
$ws.Cells.Item(4, 35).Value = 0.27
$ws.Cells.Item(4, 36).Value = 0.1
$ws.Cells.Item(4, 37).Value = 0.317
$ws.Cells.Item(4, 47).Value = 0.181
$ws.Cells.Item(4, 48).Value = 0.032
$ws.Cells.Item(4, 49).Value = 0.179
$ws.Cells.Item(4, 53).Value = 1.361
$ws.Cells.Item(4, 54).Value = 0.083
$ws.Cells.Item(4, 55).Value = 0.288
$ws.Cells.Item(4, 59).Value = 0.5
$ws.Cells.Item(4, 60).Value = 0.222
$ws.Cells.Item(4, 61).Value = 0.471
$ws.Cells.Item(4, 65).Value = 0.463
$ws.Cells.Item(4, 66).Value = 0.067
$ws.Cells.Item(4, 67).Value = 0.258
$ws.Cells.Item(4, 68).Value = 0.454
$ws.Cells.Item(4, 69).Value = 0.655
$ws.Cells.Item(5, 5).Value = 0.609
$ws.Cells.Item(5, 6).Value = 0.061
$ws.Cells.Item(5, 7).Value = 0.247
$ws.Cells.Item(5, 14).Value = 0.724
$ws.Cells.Item(5, 15).Value = 0.05
$ws.Cells.Item(5, 16).Value = 0.225
$ws.Cells.Item(5, 17).Value = 0.059
$ws.Cells.Item(5, 18).Value = 0.013
$ws.Cells.Item(5, 19).Value = 0.115
$ws.Cells.Item(5, 23).Value = 0.263
$ws.Cells.Item(5, 24).Value = 0.099
$ws.Cells.Item(5, 25).Value = 0.314
$ws.Cells.Item(5, 35).Value = 0.285
$ws.Cells.Item(5, 36).Value = 0.1
$ws.Cells.Item(5, 37).Value = 0.316
$ws.Cells.Item(5, 47).Value = 0.326
$ws.Cells.Item(5, 48).Value = 0.119
$ws.Cells.Item(5, 49).Value = 0.345
$ws.Cells.Item(5, 53).Value = 1.032
$ws.Cells.Item(5, 54).Value = 0.037
$ws.Cells.Item(5, 55).Value = 0.194
$ws.Cells.Item(5, 59).Value = 0.326
$ws.Cells.Item(5, 60).Value = 0.08799999999999999
$ws.Cells.Item(5, 61).Value = 0.297
$ws.Cells.Item(5, 65).Value = 0.365
$ws.Cells.Item(5, 66).Value = 0.036
$ws.Cells.Item(5, 67).Value = 0.189
$ws.Cells.Item(5, 68).Value = 0.344
$ws.Cells.Item(5, 69).Value = 0.395
$ws.Cells.Item(6, 5).Value = 0.583
$ws.Cells.Item(6, 14).Value = 0.55
$ws.Cells.Item(6, 17).Value = 0.08400000000000001
$ws.Cells.Item(6, 23).Value = 0.279
$ws.Cells.Item(6, 35).Value = 0.277
$ws.Cells.Item(6, 47).Value = 0.233
$ws.Cells.Item(6, 53).Value = 1.17
$ws.Cells.Item(6, 59).Value = 0.395
$ws.Cells.Item(6, 65).Value = 0.408
$ws.Cells.Item(6, 68).Value = 0.39
$ws.Cells.Item(6, 69).Value = 0.491
$ws.Cells.Item(7, 5).Value = 0.598
$ws.Cells.Item(7, 14).Value = 0.643
$ws.Cells.Item(7, 17).Value = 0.067
$ws.Cells.Item(7, 23).Value = 0.269
$ws.Cells.Item(7, 35).Value = 0.282
$ws.Cells.Item(7, 47).Value = 0.281
$ws.Cells.Item(7, 53).Value = 1.082
$ws.Cells.Item(7, 59).Value = 0.35
$ws.Cells.Item(7, 65).Value = 0.381
$ws.Cells.Item(7, 68).Value = 0.361
$ws.Cells.Item(7, 69).Value = 0.428
$ws.Cells.Item(8, 5).Value = 0.59
$ws.Cells.Item(8, 6).Value = 0.096
$ws.Cells.Item(8, 7).Value = 0.309
$ws.Cells.Item(8, 14).Value = 0.88
$ws.Cells.Item(8, 15).Value = 0.021
$ws.Cells.Item(8, 16).Value = 0.146
$ws.Cells.Item(8, 17).Value = 0.08699999999999999
$ws.Cells.Item(8, 18).Value = 0.045
$ws.Cells.Item(8, 19).Value = 0.212
$ws.Cells.Item(8, 23).Value = 0.369
$ws.Cells.Item(8, 24).Value = 0.123
$ws.Cells.Item(8, 25).Value = 0.351
$ws.Cells.Item(8, 35).Value = 0.316
$ws.Cells.Item(8, 36).Value = 0.139
$ws.Cells.Item(8, 37).Value = 0.373
$ws.Cells.Item(8, 47).Value = 0.245
$ws.Cells.Item(8, 48).Value = 0.105
$ws.Cells.Item(8, 49).Value = 0.324
$ws.Cells.Item(8, 53).Value = 1.41
$ws.Cells.Item(8, 54).Value = 0.116
$ws.Cells.Item(8, 55).Value = 0.34
$ws.Cells.Item(8, 59).Value = 0.418
$ws.Cells.Item(8, 60).Value = 0.174
$ws.Cells.Item(8, 61).Value = 0.417
$ws.Cells.Item(8, 65).Value = 0.526
$ws.Cells.Item(8, 66).Value = 0.089
$ws.Cells.Item(8, 67).Value = 0.298
$ws.Cells.Item(8, 68).Value = 0.47
$ws.Cells.Item(8, 69).Value = 0.545
$ws.Cells.Item(9, 5).Value = 0.444
$ws.Cells.Item(9, 6).Value = 0.247
$ws.Cells.Item(9, 7).Value = 0.497
$ws.Cells.Item(9, 14).Value = 0.889
$ws.Cells.Item(9, 15).Value = 0.099
$ws.Cells.Item(9, 16).Value = 0.314
$ws.Cells.Item(9, 35).Value = 0.333
$ws.Cells.Item(9, 36).Value = 0.222
$ws.Cells.Item(9, 37).Value = 0.471
$ws.Cells.Item(9, 53).Value = 1.444
$ws.Cells.Item(9, 65).Value = 0.556
$ws.Cells.Item(9, 68).Value = 0.481
$ws.Cells.Item(9, 69).Value = 0.544
$ws.Cells.Item(10, 5).Value = 0.556
$ws.Cells.Item(10, 14).Value = 1
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(10, 23).Value = 0.444
$ws.Cells.Item(10, 24).Value = 0.247
$ws.Cells.Item(10, 25).Value = 0.497
$ws.Cells.Item(10, 35).Value = 0.333
$ws.Cells.Item(10, 36).Value = 0.222
$ws.Cells.Item(10, 37).Value = 0.471
$ws.Cells.Item(10, 53).Value = 1.667
$ws.Cells.Item(10, 65).Value = 0.667
$ws.Cells.Item(10, 66).Value = 0.222
$ws.Cells.Item(10, 67).Value = 0.471
$ws.Cells.Item(10, 68).Value = 0.556
$ws.Cells.Item(10, 69).Value = 0.667
$ws.Cells.Item(11, 5).Value = 0.556
$ws.Cells.Item(11, 14).Value = 1
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = 0
$ws.Cells.Item(11, 23).Value = 0.444
$ws.Cells.Item(11, 24).Value = 0.247
$ws.Cells.Item(11, 25).Value = 0.497
$ws.Cells.Item(11, 35).Value = 0.333
$ws.Cells.Item(11, 36).Value = 0.222
$ws.Cells.Item(11, 37).Value = 0.471
$ws.Cells.Item(11, 47).Value = 0.222
$ws.Cells.Item(11, 48).Value = 0.173
$ws.Cells.Item(11, 49).Value = 0.416
$ws.Cells.Item(11, 53).Value = 1.667
$ws.Cells.Item(11, 65).Value = 0.667
$ws.Cells.Item(11, 66).Value = 0.222
$ws.Cells.Item(11, 67).Value = 0.471
$ws.Cells.Item(11, 68).Value = 0.556
$ws.Cells.Item(11, 69).Value = 0.667
$ws.Cells.Item(12, 5).Value = 1.4
$ws.Cells.Item(12, 6).Value = 0.64
$ws.Cells.Item(12, 7).Value = 0.8
$ws.Cells.Item(12, 14).Value = 1.111
$ws.Cells.Item(12, 15).Value = 0.099
$ws.Cells.Item(12, 16).Value = 0.314
$ws.Cells.Item(12, 23).Value = 1.25
$ws.Cells.Item(12, 24).Value = 0.188
$ws.Cells.Item(12, 25).Value = 0.433
$ws.Cells.Item(12, 48).Value = 6
$ws.Cells.Item(12, 49).Value = 2.449
$ws.Cells.Item(12, 53).Value = 3.367
$ws.Cells.Item(12, 54).Value = 0.16
$ws.Cells.Item(12, 55).Value = 0.4
$ws.Cells.Item(12, 65).Value = 1.167
$ws.Cells.Item(12, 66).Value = 0.139
$ws.Cells.Item(12, 67).Value = 0.373
$ws.Cells.Item(12, 68).Value = 1.122
$ws.Cells.Item(12, 69).Value = 1.231
$ws.Cells.Item(13, 68).Value = 0.803
$ws.Cells.Item(13, 69).Value = 0.68
